# Junction_Flooding_395.xlsx edit: refresh the 33-junction flooding sample with a new
# 4-row dataset (dropping the old 5th data row), "custom accuracy" column-width tweak,
# and the trimmed sheet dimension that goes with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width tweaks ("custom accuracy") -----------------------------------
# The COM ColumnWidth (character units) is stored in the OOXML <col width=..>
# attribute as ColumnWidth + 5/6 (the fixed pixel padding for this font/DPI), so we
# back that offset out to land exactly on the integer widths from the target file.
$padding = 5 / 6
$newColWidths = @{
    "B" = 8
    "C" = 8
    "F" = 7
    "G" = 8
    "I" = 8
    "K" = 8
    "L" = 8
    "M" = 8
    "O" = 8
    "P" = 8
    "T" = 9
    "V" = 8
    "W" = 8
    "X" = 8
    "AC" = 8
    "AD" = 8
    "AH" = 8
}
foreach ($col in $newColWidths.Keys) {
    $ws.Range("$col" + "1").ColumnWidth = $newColWidths[$col] - $padding
}

# --- Replace the data rows (1000-dataset refresh) -------------------------------
# Column order across each row is A (Time) then B..AH (J1..J33).
$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH")

$newRows = @{
    2 = @(45076.50694444445, 0.596, 0.5659999999999999, 0.131, 0.703, 0.594, 0, 0.673, 2.927, 1.311, 0.574, 0.838, 0.092, 0.061, 0.781, 0.09, 0.143, 2.22, 0.626, 2.165, 0.66, 0.638, 1.735, 2.278, 0.128, 0.419, 0.467, 0.214, 0.476, 0.5600000000000001, 0.8090000000000001, 2.624, 0.591, 0.381)
    3 = @(45076.51388888889, 12.809, 9.635999999999999, 0.472, 27.5, 22.92, 9.599, 31.338, 16.196, 7.646, 10.506, 11.317, 11.617, 3.215, 10.302, 13.948, 7.985, 1.254, 0.503, 146.646, 27.566, 9.500999999999999, 19.361, 10.618, 1.336, 16.444, 8.259, 7.142, 8.528, 11.835, 0.343, 28.801, 5.578, 11.438)
    4 = @(45076.52083333334, 15.248, 11.464, 0.534, 32.966, 27.4, 11.675, 45.502, 18.898, 8.795999999999999, 12.481, 13.416, 13.964, 3.846, 12.16, 16.847, 9.640000000000001, 0.917, 0.478, 175.514, 33.253, 11.226, 22.998, 12.276, 1.579, 22.665, 9.803000000000001, 8.542, 10.143, 14.08, 0.223, 41.981, 6.512, 13.689)
    5 = @(45076.52777777778, 12.39, 9.31, 0.43, 26.82, 22.29, 9.51, 39.31, 15.32, 7.13, 10.14, 10.9, 11.36, 3.13, 9.880000000000001, 13.71, 7.83, 0.71, 0.36, 141.27, 27.09, 9.130000000000001, 18.71, 9.94, 1.28, 19.11, 7.97, 6.95, 8.26, 11.45, 0.17, 36.16, 5.29, 11.13)
}

foreach ($r in $newRows.Keys) {
    $rowValues = $newRows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item([int]$r, $i + 1).Value = $rowValues[$i]
    }
}

# --- Drop the now-unused 5th data row (old row 6) -------------------------------
$ws.Rows.Item(6).Delete()

